# Rename the sheet from "Sheet1" to "Halls"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Halls"

# Move the active selection from E6 to E32
$ws.Range("E32").Select()
